$d = $word.ActiveDocument

# 1) "Megrendelő: asdasd" -> "Megrendelő: Whastz the fuck you"
$d.Content.Find.Execute("Megrendelő: asdasd", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Megrendelő: Whastz the fuck you", 2)

# 2) "Cím: " -> "Cím: Xdddddd"
$d.Content.Find.Execute("Cím: ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Cím: Xdddddd", 2)

# 3) All occurrences of "2024.06.10" -> "2024.07.01" (date updated throughout)
$d.Content.Find.Execute("2024.06.10", $false, $false, $false, $false, $false,
                         $true, 1, $false, "2024.07.01", 2)
